$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.442.14"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.832.10"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -2.98%  "
$ws.Range("D5").Value = "'315.22"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("D7").Value = "'0.4298"
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("D8").Value = "'0.3702"
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("D9").Value = "'0.07262"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "'0.8668"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "'21.18"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "1.826.97"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'5.360"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "'0.07069"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "'87.96"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").Value = "'0.000008917"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").Value = "'15.23"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Value = "27.444.93"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").Value = "'10.91"
$ws.Range("E23").Value = "  -3.23%  "
$ws.Range("D24").Value = "2.061.50"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "'2.001"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("D26").Value = "'153.48"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("D27").Value = "'18.46"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").Value = "'2.153"
$ws.Range("E28").Value = "  +7.63%  "
$ws.Range("D29").Value = "'5.289"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "'117.29"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "'0.08857"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").Value = "'1.211"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "'0.7664"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'4.488"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").Value = "'2.908"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("D37").Value = "'1.123"
$ws.Range("E37").Value = "  -2.34%  "
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "'0.05290"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'7.183"
$ws.Range("E40").Value = "  +4.42%  "
$ws.Range("D41").Value = "'2.869"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "'0.5085"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "'8.648"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "'10.57"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").Value = "'0.4747"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'106.26"
$ws.Range("E47").Value = "  -4.05%  "
$ws.Range("D48").Value = "'0.06426"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").Value = "'1.831"
$ws.Range("E51").Value = "  -3.43%  "
